$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) / Volume (E) columns hold plain text in this sheet (values
# like "225.56" or "  +2.24%  " are inline strings, not numbers/percentages).
# Cells whose new text would otherwise be auto-parsed as a number need the
# number format switched to Text first so Excel keeps storing literal text.

$ws.Range("D2").Value = "33.979.40"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.779.60"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.56"
$ws.Range("E5").Value = "  +2.24%  "
$ws.Range("E6").Value = "  +1.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.22"
$ws.Range("E8").Value = "  +4.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.291"
$ws.Range("E9").Value = "  +2.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0701"
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("D12").Value = "2.036.94"
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.05"
$ws.Range("E13").Value = "  +5.51%  "
$ws.Range("D14").Value = "1.777.29"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.621"
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").Value = "33.972.92"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.16"
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.62"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.41"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "0.0₃0782"
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.70"
$ws.Range("E22").Value = "  +1.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.11"
$ws.Range("E23").Value = "  +2.40%  "
$ws.Range("E24").Value = "  -2.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.56"
$ws.Range("E25").Value = "  +1.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.26"
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("E27").Value = "  +1.86%  "
$ws.Range("E28").Value = "  +1.27%  "
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.23"
$ws.Range("E30").Value = "  +3.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0512"
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("E32").Value = "  -0.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.56"
$ws.Range("E33").Value = "  +1.98%  "
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").Value = "1.393.22"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.661"
$ws.Range("E36").Value = "  +5.35%  "
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("E38").Value = "  +0.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.24"
$ws.Range("E39").Value = "  +6.90%  "
$ws.Range("E40").Value = "  +1.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.912"
$ws.Range("E41").Value = "  -1.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "78.02"
$ws.Range("E42").Value = "  -0.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.65"
$ws.Range("E43").Value = "  -1.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.47"
$ws.Range("E44").Value = "  +14.06%  "
$ws.Range("D45").Value = "0.0₆0146"
$ws.Range("E45").Value = "  +22.61%  "
$ws.Range("E46").Value = "  +4.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "108.56"
$ws.Range("E47").Value = "  +4.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0498"
$ws.Range("E48").Value = "  +1.35%  "
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("D50").Value = "1.936.61"
$ws.Range("E50").Value = "  +1.00%  "
$ws.Range("E51").Value = "  +0.51%  "
